# "Add files via upload" — re-saved workbook with updated sample data in
# the block of cells around M11:Q15, plus the active selection left on N14
# (where the author was last working) before upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (M..Q): 0, 80, 80, 80, 0  ->  0, 40, 40, 40, 0
$ws.Range("N11:P11").Value = 40

# Row 12 (M..Q): 80, 80, 235, 80, 80  ->  40, 40, 235.9, 40, 40
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 40
$ws.Range("O12").Value = 235.9
$ws.Range("P12").Value = 40
$ws.Range("Q12").Value = 40

# Row 13 (M..Q): 80, 235, 235, 235, 80  ->  40, 235.9, 235.9, 235.9, 40
$ws.Range("M13").Value = 40
$ws.Range("N13").Value = 235.9
$ws.Range("O13").Value = 235.9
$ws.Range("P13").Value = 235.9
$ws.Range("Q13").Value = 40

# Row 14 (M..Q): 80, 80, 235, 80, 80  ->  40, 40, 235.9, 40, 40
$ws.Range("M14").Value = 40
$ws.Range("N14").Value = 40
$ws.Range("O14").Value = 235.9
$ws.Range("P14").Value = 40
$ws.Range("Q14").Value = 40

# Row 15 (M..Q): 0, 80, 80, 80, 0  ->  0, 40, 40, 40, 0
$ws.Range("N15:P15").Value = 40

# Leave the active cell/selection on N14, matching the saved view state.
$ws.Range("N14").Select()
